# Generate Report for Handoff
# The localization status moved from "In Translation" to "Ready for handoff"
# and the handoff timestamps were refreshed. Update the three report
# sheets (Overview, zh-cn, de-de) accordingly and re-fit the now-wider
# "Status" columns.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet --------------------------------------------------
# E2/F2 hold the per-language status, G2 the latest handoff xliff
# generation timestamp.
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-15 22:54:52"

# --- zh-cn detail sheet ------------------------------------------------
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-15 22:54:48"

# --- de-de detail sheet ------------------------------------------------
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-15 22:54:52"

# --- Widen the Status columns to fit "Ready for handoff" ---------------
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333
$zhcn.Columns.Item(3).ColumnWidth = 16.3333333333333
$dede.Columns.Item(3).ColumnWidth = 16.3333333333333
